# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# Rename header labels on existing sheets
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add($null, $lastSheet)
$new.Name = "PO Forecast"

# Copy header/date formatting from the Weekly Quantity sheet so styles line up
$ws1.Range("A1:B1").Copy()
$new.Range("A1:D1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$new.Range("A2:A18").PasteSpecial(-4122)

# Header labels
$new.Range("A1").Value = "ds"
$new.Range("B1").Value = "PO_Forecast"
$new.Range("C1").Value = "yhat_lower"
$new.Range("D1").Value = "yhat_upper"

$new.Range("A2").Value = 44941.99999999999
$new.Range("B2").Value = 332
$new.Range("C2").Value = 253.4522132684213
$new.Range("D2").Value = 411.3155457228025
$new.Range("A3").Value = 44962.99999999999
$new.Range("B3").Value = 252
$new.Range("C3").Value = 170.2949249094218
$new.Range("D3").Value = 330.0871512403432
$new.Range("A4").Value = 44976.99999999999
$new.Range("B4").Value = 199
$new.Range("C4").Value = 121.7718260428648
$new.Range("D4").Value = 273.9830998178775
$new.Range("A5").Value = 44990.99999999999
$new.Range("B5").Value = 145
$new.Range("C5").Value = 68.66528966001948
$new.Range("D5").Value = 223.8180409423034
$new.Range("A6").Value = 44997.99999999999
$new.Range("B6").Value = 119
$new.Range("C6").Value = 42.6595592112245
$new.Range("D6").Value = 195.5153102351869
$new.Range("A7").Value = 45004.99999999999
$new.Range("B7").Value = 92
$new.Range("C7").Value = 15.05236162447617
$new.Range("D7").Value = 170.3544847003631
$new.Range("A8").Value = 45011.99999999999
$new.Range("B8").Value = 65
$new.Range("C8").Value = -9.918576947272994
$new.Range("D8").Value = 141.9666229086928
$new.Range("A9").Value = 45018.99999999999
$new.Range("B9").Value = 39
$new.Range("C9").Value = -41.1820848932865
$new.Range("D9").Value = 116.9196120324288
$new.Range("A10").Value = 45032.99999999999
$new.Range("B10").Value = 0
$new.Range("C10").Value = -92.58177090536009
$new.Range("D10").Value = 59.41036703003063
$new.Range("A11").Value = 45039.99999999999
$new.Range("B11").Value = 0
$new.Range("C11").Value = -121.5975907961046
$new.Range("D11").Value = 39.19518382188178
$new.Range("A12").Value = 45046.99999999999
$new.Range("B12").Value = 0
$new.Range("C12").Value = -143.5349575903837
$new.Range("D12").Value = 7.667615098059104
$new.Range("A13").Value = 45053.99999999999
$new.Range("B13").Value = 0
$new.Range("C13").Value = -172.4864538994149
$new.Range("D13").Value = -18.80595458126125
$new.Range("A14").Value = 45060.99999999999
$new.Range("B14").Value = 0
$new.Range("C14").Value = -200.2625654441687
$new.Range("D14").Value = -46.29101847968764
$new.Range("A15").Value = 45067.99999999999
$new.Range("B15").Value = 0
$new.Range("C15").Value = -232.9277507735973
$new.Range("D15").Value = -70.15271535148342
$new.Range("A16").Value = 45074.99999999999
$new.Range("B16").Value = 0
$new.Range("C16").Value = -250.5120089659029
$new.Range("D16").Value = -96.71608761770345
$new.Range("A17").Value = 45081.99999999999
$new.Range("B17").Value = 0
$new.Range("C17").Value = -281.8327325502055
$new.Range("D17").Value = -129.0203774429511
$new.Range("A18").Value = 45088.99999999999
$new.Range("B18").Value = 0
$new.Range("C18").Value = -307.2806077507333
$new.Range("D18").Value = -152.9657283623877

# Restore the originally active sheet/selection
$ws1.Activate()
